$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "OBJECTIVE:" paragraph right before "EDUCATION:"
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("EDUCATION:") | Out-Null
$eduPara = $findRange.Paragraphs.First
$eduPara.Range.InsertParagraphBefore()

$allParas = $d.Paragraphs
$objPara = $allParas.Item(4)
$objPara.Style = "Heading1"

$objXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">OBJECTIVE: </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Eager to drive back-end solutions at Zion Bank on a full-time basis</w:t></w:r>' +
  '</w:p></w:body></w:document>'

$objPara.Range.InsertXML($objXml)

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the document to the
#    middle of "Computer Science" (".S. in Computer Scie|nce")
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$csRange = $d.Content
$csRange.Find.Execute("Computer Science") | Out-Null
$splitPos = $csRange.Start + "Computer Scie".Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
